$d = $word.ActiveDocument

# The document has three inline pictures living in the headers/footers:
#   - Footer 1: Pearson Edexcel logo (docPr id=1) -> rename "image1.png" to "image2.png"
#   - Footer 2: Pearson Edexcel logo (docPr id=2) -> rename "image1.png" to "image2.png"
#   - Header 2: BTec logo            (docPr id=3) -> rename "image2.jpg" to "image1.jpg"
#
# Walk every section's headers/footers and rename the inline picture(s)
# found there based on their current (old) name, so the script is robust
# even if headers/footers/sections are enumerated in a different order.

function Rename-InlineShapesInRange($range, $oldName, $newName) {
    for ($i = 1; $i -le $range.InlineShapes.Count; $i++) {
        $shp = $range.InlineShapes.Item($i)
        if ($shp.Name -eq $oldName -or [string]::IsNullOrEmpty($shp.Name)) {
            $shp.Name = $newName
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
                $shp = $hdr.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    for ($f = 1; $f -le 3; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
                $shp = $ftr.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}

# Also cover any inline pictures living directly in the main document body,
# in case layout differs from what was inspected.
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
        $shp.Name = "image1.jpg"
    } elseif ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shp.Name = "image2.png"
    }
}
